# AccessModifyingProperties.pptx - remove Aspose.Slides evaluation
# watermarks that were left over on the two slides, and replace the
# watermark shape on slide 2 with an (empty) placeholder textbox, as
# produced when the "Evaluation only." trial-version stamp is removed.

$p = $ppt.ActivePresentation

# --- Slide 1: delete the "Evaluation only." watermark textbox (id 2055) ---
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s1.Shapes.Item($i)
    if ($sh.Name -eq "TextBox" -and $sh.Id -eq 2055) {
        $sh.Delete()
    }
}

# --- Slide 2: clear + reposition the watermark textbox, add a new blank textbox ---
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.Name -eq "TextBox" -and $sh.Id -eq 3076) {
        $sh.TextFrame.TextRange.Text = ""
        $sh.Left = 4479841 / 12700.0
        $sh.Top = 3051623 / 12700.0
        $sh.Width = 184730 / 12700.0
        $sh.Height = 754694 / 12700.0
    }
}

$nb = $s2.Shapes.AddTextbox(1, 4427984 / 12700.0, 2780928 / 12700.0, 184731 / 12700.0, 369332 / 12700.0)
$nb.Fill.Visible = 0
$nb.TextFrame.WordWrap = 0
$nb.TextFrame.AutoSize = 1
